$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("Q3")
$c.Value = 2020
$c.Style = "Обычный"
$c.Font.Name = "Times New Roman"
Write-Output "done"
